$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 119, shifting existing rows 119-171 down to 120-172.
$ws.Rows.Item(119).Insert()

# Populate the newly-inserted row 119 with the new weekly record.
# Columns that stay identical to the (now shifted-down) former row 119 / row 120
# are re-filled explicitly since Insert() leaves the new row blank.
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44523
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = 100112021
$ws.Range("G119").Value = "Ají"
$ws.Range("H119").Value = "Inferno"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 140
$ws.Range("K119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("M119").Value = 25000
$ws.Range("N119").Value = "$/caja 12 kilos"
$ws.Range("O119").Value = "Región de Arica y Parinacota"
$ws.Range("P119").Value = 2083
$ws.Range("Q119").Value = 12
$ws.Range("R119").Value = "Hortaliza"
